# issue #5: add legislator_id, name, date into dataframe
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("股票")

# Header row: new columns H (date), I (legislator_name), J (legislator_id)
$ws.Cells.Item(1, 8).Value = "date"
$ws.Cells.Item(1, 9).Value = "legislator_name"
$ws.Cells.Item(1, 10).Value = "legislator_id"

# Force column H to be treated as plain text so "2012-03-30" isn't
# auto-converted into a date serial number.
$dateRange = $ws.Range("H2:H8")
$dateRange.NumberFormat = "@"

for ($r = 2; $r -le 8; $r++) {
    $ws.Cells.Item($r, 8).Value = "2012-03-30"
    $ws.Cells.Item($r, 9).Value = "薛凌"
    $ws.Cells.Item($r, 10).Value = 1384
}

# Drop the temporary text format so the cells fall back to the sheet's
# normal (default) style instead of keeping a bespoke "@" number format.
$dateRange.ClearFormats()
